$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.455.55'
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.870.75'
$ws.Range("E3").Value = '  -0.47%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7089'
$ws.Range("E5").Value = '  -0.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.87'
$ws.Range("E6").Value = '  +0.56%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3166'
$ws.Range("E8").Value = '  +0.85%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07887'
$ws.Range("E9").Value = '  -1.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.69'
$ws.Range("E10").Value = '  -1.73%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08003'
$ws.Range("E11").Value = '  -3.85%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.887.44'
$ws.Range("E12").Value = '  +0.00%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.232'
$ws.Range("E13").Value = '  -0.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.16'
$ws.Range("E14").Value = '  -0.32%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7063'
$ws.Range("E15").Value = '  -1.63%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.524'
$ws.Range("E16").Value = '  +2.60%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.474.53'
$ws.Range("E17").Value = '  +0.40%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008367'
$ws.Range("E18").Value = '  -4.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '257.31'
$ws.Range("E19").Value = '  +5.93%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.128.53'
$ws.Range("E20").Value = '  -1.04%  '

$ws.Range("E21").Value = '  -0.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.640'
$ws.Range("E23").Value = '  -2.70%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.0000'
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1563'
$ws.Range("E25").Value = '  -0.50%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.086'
$ws.Range("E26").Value = '  +0.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.98'
$ws.Range("E27").Value = '  -1.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.93'
$ws.Range("E28").Value = '  +1.80%  '

$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.346'
$ws.Range("E30").Value = '  -2.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.264'
$ws.Range("E31").Value = '  -2.16%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.211'
$ws.Range("E32").Value = '  +0.56%  '

$ws.Range("E33").Value = '  -1.26%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.901'
$ws.Range("E34").Value = '  -2.19%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7515'
$ws.Range("E35").Value = '  -3.38%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.177'
$ws.Range("E36").Value = '  -0.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.714'
$ws.Range("E37").Value = '  +0.93%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01886'
$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.273.62'
$ws.Range("E39").Value = '  +0.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.755'
$ws.Range("E40").Value = '  +0.36%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9019'
$ws.Range("E41").Value = '  -1.93%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.987'
$ws.Range("E42").Value = '  -8.62%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '108.23'
$ws.Range("E43").Value = '  -4.96%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.76'
$ws.Range("E44").Value = '  -3.80%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9996'
$ws.Range("E45").Value = '  -0.14%  '

$ws.Range("E46").Value = '  +1.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.027.48'
$ws.Range("E47").Value = '  -0.47%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.796'
$ws.Range("E48").Value = '  -0.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5195'
$ws.Range("E49").Value = '  -0.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.553'
$ws.Range("E50").Value = '  -0.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4338'
$ws.Range("E51").Value = '  -1.00%  '

Write-Host "Applied cryptos update"